$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These cells hold numeric-looking values that are stored as TEXT
# (shared strings) in the workbook. When assigning a plain numeric
# string to Range.Value, Excel auto-converts it to a real number,
# which would change the underlying cell type. To preserve the
# original text type we:
#   1. remember the cell's current Style
#   2. assign the value with a leading apostrophe so Excel keeps it
#      as text (quote-prefixed)
#   3. restore the original Style so formatting/appearance is
#      unchanged

function Set-TextValue {
    param($cell, [string]$text)
    $origStyle = $cell.Style
    $cell.Value = "'" + $text
    $cell.Style = $origStyle
}

# "Enterprises density (per 1000 people)" row (Statistical Institution table)
Set-TextValue $ws.Range("B11") "30.94"
Set-TextValue $ws.Range("C11") "1.24"
Set-TextValue $ws.Range("D11") "32.18"

# "Enterprises (% of total)" row (Statistical Institution table)
Set-TextValue $ws.Range("B13") "96.01"
Set-TextValue $ws.Range("C13") "3.86"
Set-TextValue $ws.Range("D13") "99.87"

# "Enterprises density (per 1000 people)" row (SME Associations table)
Set-TextValue $ws.Range("C32") "0.62"
Set-TextValue $ws.Range("D32") "14.42"
